$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shp = $s.Shapes.Item("Content Placeholder 4")
$tr = $shp.TextFrame.TextRange

# Paragraph 3 currently reads: "Scala Worksheets in Scala IDE"
# It should become: "Scala Worksheets in Intelijj. Set Scala SDK to Scala 2.11"
$para = $tr.Paragraphs(3, 1)

# Shrink the trailing "in Scala IDE" run down to just "in " (keep its own run/formatting).
$tail = $para.Characters(18, 13)
$tail.Text = "in "

# Append the remaining new text as its own separate runs, re-fetching the
# paragraph reference each time so every InsertAfter() lands on its own run.
$para = $tr.Paragraphs(3, 1)
[void]$para.InsertAfter("Intelijj")

$para = $tr.Paragraphs(3, 1)
[void]$para.InsertAfter(". Set Scala SDK ")

$para = $tr.Paragraphs(3, 1)
[void]$para.InsertAfter("to Scala 2.11")
